$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap columns B and D (header + all data rows)
$tmp = $ws.Range("B1:B7").Value2
$ws.Range("B1:B7").Value2 = $ws.Range("D1:D7").Value2
$ws.Range("D1:D7").Value2 = $tmp

# Swap columns E and F (header + all data rows)
$tmp2 = $ws.Range("E1:E7").Value2
$ws.Range("E1:E7").Value2 = $ws.Range("F1:F7").Value2
$ws.Range("F1:F7").Value2 = $tmp2
